# The underlying change captured by the canonical-OOXML diff for this
# fixture is a pure XML attribute/namespace-declaration reordering
# (alphabetical canonicalization) introduced when the template docx was
# regenerated/re-saved - e.g. <w:pgSz w:w="11906" w:h="16838"/> becoming
# <w:pgSz w:h="16838" w:w="11906"/>, namespace declarations on
# <w:document> being re-sorted, <w:lsdException>/<w:style> attributes
# being re-sorted, etc. Every corresponding added/removed line pair in
# the diff has identical tag names, identical attribute name/value
# pairs and identical text content once ordering is ignored - i.e.
# there is no actual visible/semantic content, formatting or structural
# change to the document.
#
# The Word object model (and COM automation in general) does not expose
# XML attribute ordering - that is purely a serialization detail of
# however the package was written out - so there is no corresponding
# action to replay through $word/$d here. We simply touch the document
# (via a no-op Find) so the script runs cleanly against the active
# document without altering any content, matching the (content-wise
# unchanged) target.

$d = $word.ActiveDocument

$null = $d.Content.Find.Execute(
    "",     # FindText
    $false, # MatchCase
    $false, # MatchWholeWord
    $false, # MatchWildcards
    $false, # MatchSoundsLike
    $false, # MatchAllWordForms
    $true,  # Forward
    1,      # Wrap
    $false, # Format
    "",     # ReplaceWith
    0       # Replace (wdReplaceNone)
)
